$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.901.56"
$ws.Range("E2").Value = "  +1.69%  "

$ws.Range("D3").Value = "1.676.29"
$ws.Range("E3").Value = "  +1.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9943"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3647"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07089"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9993"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.101"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.88%  "

$ws.Range("D15").Value = "1.676.02"
$ws.Range("E15").Value = "  +1.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.631"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001054"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06601"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9967"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "79.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.28%  "

$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.53%  "

$ws.Range("D24").Value = "24.965.35"
$ws.Range("E24").Value = "  +2.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.449"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.426"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.77%  "

$ws.Range("D29").Value = "1.860.04"
$ws.Range("E29").Value = "  +1.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.195"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.060"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.787"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08455"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.644"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.184"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02272"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.53%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06075"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.233"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.32%  "

$ws.Range("E41").Value = "  +2.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.264"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9941"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5974"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.833"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.95%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5724"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.963"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07031"
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.191"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.52%  "
